$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "surcharge" column header to "network_surcharge"
$ws.Range("H1").Value = "network_surcharge"

# Update the active selection to just the H1 cell
$ws.Range("H1").Select()
